$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on the Price/Volume columns so that numeric-
# looking strings (e.g. "1.001") are stored as text, matching the original
# inlineStr cell type, instead of being auto-converted to numbers by Excel.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '30.557.33'
$ws.Range("E2").Value = '  -1.03%  '
$ws.Range("D3").Value = '1.875.84'
$ws.Range("E3").Value = '  -1.38%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '235.54'
$ws.Range("E5").Value = '  -4.42%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("D7").Value = '0.4867'
$ws.Range("E7").Value = '  -2.54%  '
$ws.Range("D8").Value = '0.2902'
$ws.Range("E8").Value = '  -2.85%  '
$ws.Range("E9").Value = '  -2.82%  '
$ws.Range("D10").Value = '1.876.74'
$ws.Range("E10").Value = '  -1.44%  '
$ws.Range("D11").Value = '16.56'
$ws.Range("E11").Value = '  -4.22%  '
$ws.Range("E12").Value = '  -1.52%  '
$ws.Range("D13").Value = '88.86'
$ws.Range("E13").Value = '  -3.47%  '
$ws.Range("D14").Value = '4.988'
$ws.Range("D15").Value = '0.6520'
$ws.Range("E15").Value = '  -4.61%  '
$ws.Range("D16").Value = '30.502.07'
$ws.Range("E16").Value = '  -1.27%  '
$ws.Range("D17").Value = '0.000007846'
$ws.Range("E17").Value = '  -3.05%  '
$ws.Range("D18").Value = '1.001'
$ws.Range("E18").Value = '  +0.17%  '
$ws.Range("D19").Value = '12.92'
$ws.Range("E19").Value = '  -3.93%  '
$ws.Range("D20").Value = '2.121.68'
$ws.Range("E20").Value = '  -1.43%  '
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  +0.15%  '
$ws.Range("D22").Value = '4.700'
$ws.Range("E22").Value = '  -3.80%  '
$ws.Range("D23").Value = '191.37'
$ws.Range("E23").Value = '  +5.22%  '
$ws.Range("D24").Value = '6.080'
$ws.Range("E24").Value = '  -0.56%  '
$ws.Range("D25").Value = '9.255'
$ws.Range("E25").Value = '  -1.32%  '
$ws.Range("D26").Value = '158.15'
$ws.Range("E26").Value = '  +2.09%  '
$ws.Range("D27").Value = '18.29'
$ws.Range("E27").Value = '  -2.62%  '
$ws.Range("D28").Value = '1.819'
$ws.Range("E28").Value = '  -6.82%  '
$ws.Range("D29").Value = '1.403'
$ws.Range("E29").Value = '  +0.57%  '
$ws.Range("D30").Value = '4.228'
$ws.Range("E30").Value = '  -3.57%  '
$ws.Range("D31").Value = '0.08995'
$ws.Range("E31").Value = '  +0.18%  '
$ws.Range("D32").Value = '3.914'
$ws.Range("E32").Value = '  -3.82%  '
$ws.Range("D33").Value = '0.05128'
$ws.Range("E33").Value = '  -3.55%  '
$ws.Range("D34").Value = '0.7220'
$ws.Range("E34").Value = '  -4.08%  '
$ws.Range("E35").Value = '  -6.07%  '
$ws.Range("D36").Value = '2.694'
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("E37").Value = '  -5.88%  '
$ws.Range("D38").Value = '2.655'
$ws.Range("E38").Value = '  -2.73%  '
$ws.Range("D39").Value = '0.9180'
$ws.Range("E39").Value = '  -2.39%  '
$ws.Range("D40").Value = '2.041'
$ws.Range("E40").Value = '  -7.20%  '
$ws.Range("D41").Value = '0.4364'
$ws.Range("E41").Value = '  -0.85%  '
$ws.Range("D42").Value = '104.45'
$ws.Range("E42").Value = '  -1.73%  '
$ws.Range("D43").Value = '0.9950'
$ws.Range("E43").Value = '  -0.55%  '
$ws.Range("D44").Value = '5.693'
$ws.Range("E44").Value = '  -3.34%  '
$ws.Range("D45").Value = '0.1326'
$ws.Range("E45").Value = '  -4.13%  '
$ws.Range("D46").Value = '7.322'
$ws.Range("E46").Value = '  -5.83%  '
$ws.Range("D47").Value = '0.4027'
$ws.Range("E47").Value = '  +2.23%  '
$ws.Range("D48").Value = '0.05820'
$ws.Range("E48").Value = '  -0.46%  '
$ws.Range("D49").Value = '8.650'
$ws.Range("E49").Value = '  +0.53%  '
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("D51").Value = '33.08'
$ws.Range("E51").Value = '  -1.58%  '

# Restore the default cell style so no stray formatting is introduced.
$dataRange.Style = "Normal"

